$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (for line7 and line8) right after the existing
# "line6" row (row 7) and before the "extr1" row (current row 8).
# This pushes the existing extr1..extr8 rows down by two (rows 8-15 -> 10-17).
$ws.Rows("8:9").Insert()

# Copy the formatting (style) of the column-A cells from nearby rows so the
# new rows keep the same bold/border/centered style (s="1") used by all the
# other data rows in column A.
$ws.Cells.Item(10,1).Copy($ws.Cells.Item(8,1))
$ws.Cells.Item(10,1).Copy($ws.Cells.Item(9,1))

# Fill in the data for the new "line7" / "line8" rows.
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "line7"
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = 11
$ws.Cells.Item(8,5).Value = $true

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "line8"
$ws.Cells.Item(9,3).Value = 16
$ws.Cells.Item(9,4).Value = 9
$ws.Cells.Item(9,5).Value = $true

# Update the data for the (now shifted) extr1..extr8 rows (rows 10-17).
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "extr1"
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 12
$ws.Cells.Item(10,5).Value = $true

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "extr2"
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 9
$ws.Cells.Item(11,5).Value = $true

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "extr3"
$ws.Cells.Item(12,3).Value = 10
$ws.Cells.Item(12,4).Value = 11
$ws.Cells.Item(12,5).Value = $true

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "extr4"
$ws.Cells.Item(13,3).Value = 7
$ws.Cells.Item(13,4).Value = 8
$ws.Cells.Item(13,5).Value = $false

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "extr5"
$ws.Cells.Item(14,3).Value = 9
$ws.Cells.Item(14,4).Value = 11
$ws.Cells.Item(14,5).Value = $true

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "extr6"
$ws.Cells.Item(15,3).Value = 7
$ws.Cells.Item(15,4).Value = 11
$ws.Cells.Item(15,5).Value = $false

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "extr7"
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 7
$ws.Cells.Item(16,5).Value = $true

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "extr8"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = $true
